# The canonical diff shows a new record being inserted right before the
# existing row 228, pushing every subsequent row (228..310) down by one
# (to 229..311) and extending the sheet dimension to A1:R311.
#
# Replicate that with a real row insertion (which shifts the rows below
# down and carries their formatting/styles along, exactly like Excel's
# own Insert does), then populate the newly inserted row with the data
# for the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 228; rows 228-310 shift down to 229-311.
$ws.Rows.Item(228).Insert()

# Populate the newly inserted row 228 with the new record's data.
$ws.Range("A228").Value = 4
$ws.Range("B228").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C228").Value = "Los Lagos"
$ws.Range("D228").Value = 44917
$ws.Range("E228").Value = 10
$ws.Range("F228").Value = 100112028
$ws.Range("G228").Value = "Sandia"
$ws.Range("H228").Value = "Sin especificar"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 300
$ws.Range("K228").Value = 4000
$ws.Range("L228").Value = 4000
$ws.Range("M228").Value = 4000
$ws.Range("N228").Value = "`$/unidad"
$ws.Range("O228").Value = "Región de O'Higgins"
$ws.Range("P228").Value = 4000
$ws.Range("Q228").Value = 1
$ws.Range("R228").Value = "Hortaliza"
